$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 142
$ws.Range("A142").Value = 45971
$ws.Range("B142").Value = "四方坪站充电量(kw)"
$ws.Range("C142").Value = 664.46300000000008
$ws.Range("D142").Value = 1537.576
$ws.Range("E142").Value = 548.98
$ws.Range("F142").Value = 448.92899999999997
$ws.Range("G142").Value = 231.18400000000003
$ws.Range("H142").Value = 649.60699999999997
$ws.Range("I142").Value = 611.38300000000004
$ws.Range("J142").Value = 195.26499999999999
$ws.Range("K142").Value = 157.02700000000002
$ws.Range("L142").Value = 251.97999999999996
$ws.Range("M142").Value = 229.76500000000004
$ws.Range("N142").Value = 221.51599999999996
$ws.Range("O142").Value = 936.72900000000016
$ws.Range("P142").Value = 1657.9210000000003
$ws.Range("Q142").Value = 696.02999999999986
$ws.Range("R142").Value = 558.32199999999989
$ws.Range("S142").Value = 445.42399999999998
$ws.Range("T142").Value = 198.34200000000001
$ws.Range("U142").Value = 258.58
$ws.Range("V142").Value = 54.53
$ws.Range("W142").Value = 131.39099999999999
$ws.Range("X142").Value = 73.349999999999994
$ws.Range("Y142").Value = 8.5399999999999991
$ws.Range("Z142").Value = 31.32

# Row 143
$ws.Range("A143").Value = 45971
$ws.Range("B143").Value = "高岭站充电量(kw)"
$ws.Range("C143").Value = 269.74399999999997
$ws.Range("D143").Value = 286.83199999999999
$ws.Range("E143").Value = 123.21099999999998
$ws.Range("F143").Value = 29.231999999999999
$ws.Range("G143").Value = 30.695
$ws.Range("H143").Value = 20.062999999999999
$ws.Range("I143").Value = 356.65999999999997
$ws.Range("J143").Value = 139.41399999999999
$ws.Range("K143").Value = 24.742000000000001
$ws.Range("L143").Value = 127.43600000000001
$ws.Range("M143").Value = 174.83799999999999
$ws.Range("N143").Value = 220.428
$ws.Range("O143").Value = 237.071
$ws.Range("P143").Value = 890.33700000000022
$ws.Range("Q143").Value = 297.39100000000002
$ws.Range("R143").Value = 290.61699999999996
$ws.Range("S143").Value = 371.60300000000001
$ws.Range("T143").Value = 232.33699999999999
$ws.Range("U143").Value = 64.076999999999998
$ws.Range("V143").Value = 130.917
$ws.Range("W143").Value = 29.67
$ws.Range("X143").Value = 119.809
$ws.Range("Y143").Value = 54.661999999999999
$ws.Range("Z143").Value = 33.094000000000001

$ws.Range("F151").Select()